$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-preserving cells to Text format first so Excel does not
# reinterpret numeric-looking literals (trailing zeros, etc.) as numbers.
$textCells = @('D5', 'D6', 'D14', 'D19', 'D20', 'D21', 'D23', 'D24', 'D25', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D39', 'D40', 'D42', 'D46')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (matches the upstream data refresh).
$ws.Range('D2').Value = '60.318.22'
$ws.Range('E2').Value = '  -5.04%  '
$ws.Range('D3').Value = '3.011.66'
$ws.Range('E3').Value = '  -5.11%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '573.19'
$ws.Range('E5').Value = '  -3.38%  '
$ws.Range('D6').Value = '126.11'
$ws.Range('E6').Value = '  -7.06%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.007.30'
$ws.Range('E8').Value = '  -5.14%  '
$ws.Range('E9').Value = '  -2.16%  '
$ws.Range('E10').Value = '  -7.68%  '
$ws.Range('E11').Value = '  -5.35%  '
$ws.Range('E12').Value = '  -2.85%  '
$ws.Range('E13').Value = '  -7.50%  '
$ws.Range('D14').Value = '32.68'
$ws.Range('E14').Value = '  -5.81%  '
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').Value = '3.504.50'
$ws.Range('E16').Value = '  -5.25%  '
$ws.Range('D17').Value = '3.009.91'
$ws.Range('E17').Value = '  -5.23%  '
$ws.Range('D18').Value = '60.271.30'
$ws.Range('E18').Value = '  -5.11%  '
$ws.Range('D19').Value = '6.53'
$ws.Range('E19').Value = '  -0.32%  '
$ws.Range('D20').Value = '430.54'
$ws.Range('E20').Value = '  -6.80%  '
$ws.Range('D21').Value = '13.20'
$ws.Range('E21').Value = '  -5.47%  '
$ws.Range('E22').Value = '  -3.57%  '
$ws.Range('D23').Value = '7.09'
$ws.Range('E23').Value = '  -7.82%  '
$ws.Range('D24').Value = '12.97'
$ws.Range('E24').Value = '  -2.08%  '
$ws.Range('D25').Value = '79.48'
$ws.Range('E25').Value = '  -4.38%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('E28').Value = '  -4.85%  '
$ws.Range('E29').Value = '  -3.94%  '
$ws.Range('D30').Value = '7.27'
$ws.Range('E30').Value = '  -6.24%  '
$ws.Range('D31').Value = '6.16'
$ws.Range('E31').Value = '  -10.00%  '
$ws.Range('D32').Value = '25.39'
$ws.Range('E32').Value = '  -7.22%  '
$ws.Range('D33').Value = '0.0950'
$ws.Range('E33').Value = '  -5.73%  '
$ws.Range('D34').Value = '5.63'
$ws.Range('E34').Value = '  -4.44%  '
$ws.Range('D35').Value = '0.940'
$ws.Range('E35').Value = '  -8.11%  '
$ws.Range('E37').Value = '  -15.12%  '
$ws.Range('D38').Value = '0.0₃0671'
$ws.Range('E38').Value = '  -8.21%  '
$ws.Range('D39').Value = '8.52'
$ws.Range('E39').Value = '  +4.79%  '
$ws.Range('D40').Value = '0.0357'
$ws.Range('E40').Value = '  -8.62%  '
$ws.Range('E41').Value = '  -4.19%  '
$ws.Range('D42').Value = '372.75'
$ws.Range('E42').Value = '  -5.19%  '
$ws.Range('D43').Value = '2.682.72'
$ws.Range('E43').Value = '  -3.83%  '
$ws.Range('E44').Value = '  -7.56%  '
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').Value = '121.59'
$ws.Range('E46').Value = '  -4.66%  '
$ws.Range('E47').Value = '  -6.40%  '
$ws.Range('E48').Value = '  -5.15%  '
$ws.Range('E49').Value = '  -3.31%  '
$ws.Range('E50').Value = '  -6.93%  '
$ws.Range('E51').Value = '  -6.53%  '

# Restore the original (default) formatting/style on the text-forced cells
# so only the cell values change, not their styles.
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Style = "Normal"
}